# Add new column F ("SE+AR") to the data moments table.
# FIRE moved before the stylized-fact column; a new SE+AR column is appended
# with a header and per-row moment formulas/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column, styled like the existing header cells (B1:E1)
$ws.Range("F1").Value = "SE+AR"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row values for the new column F
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = "`$\sigma^2/(1-\rho^2)`$"
$ws.Range("F4").Value = "`$\rho\sigma^2/(1-\rho^2)`$"
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = "`$\lambda^2\sigma^2/(1-(1-\lambda)^2\rho^2)`$"
$ws.Range("F7").Value = "`$(1-\lambda)\rho\text{FEVar}`$"
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = "`$\sigma^2`$"
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
